# Cập nhật tên bài tập
# Slide 23 title: "Bài tập"    -> "Bài tập 6.1"
# Slide 24 title: "Bài 2"      -> "Bài tập 6.2"

$p = $ppt.ActivePresentation

# --- Slide 23 ("Bài tập" -> "Bài tập 6.1") -------------------------------
$slide23 = $p.Slides.Item(23)
$title23 = $slide23.Shapes.Item(1)
$tr23 = $title23.TextFrame.TextRange

# Locate the run that holds "tập" (the word after "Bài ") and clear it,
# then type the new wording after it so the run boundaries / formatting
# of the untouched runs ("Bài", " ") stay intact.
$word23 = $tr23.Characters(5, 3)
$word23.Text = ""
$tr23 = $title23.TextFrame.TextRange
$tail23 = $tr23.Characters($tr23.Length, 0)
$tail23.InsertAfter("tập 6.1") | Out-Null

# --- Slide 24 ("Bài 2" -> "Bài tập 6.2") ---------------------------------
$slide24 = $p.Slides.Item(24)
$title24 = $slide24.Shapes.Item(1)
$tr24 = $title24.TextFrame.TextRange

# Remove the trailing "2" (keep the leading space that follows "Bài"),
# then append the new "tập 6.2" text as its own run.
$num24 = $tr24.Characters(5, 1)
$num24.Text = ""
$tr24 = $title24.TextFrame.TextRange
$tail24 = $tr24.Characters($tr24.Length, 0)
$tail24.InsertAfter("tập 6.2") | Out-Null
